# Update time_taken (column F) timestamps on the "data" sheet.
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value = "2021-10-05 14:35:33.681634"
$dataSheet.Range("F3").Value = "2021-10-05 14:35:33.681642"
$dataSheet.Range("F4").Value = "2021-10-05 14:35:33.681645"
$dataSheet.Range("F5").Value = "2021-10-05 14:35:33.681648"
$dataSheet.Range("F6").Value = "2021-10-05 14:35:33.681650"
$dataSheet.Range("F7").Value = "2021-10-05 14:35:33.681653"
$dataSheet.Range("F8").Value = "2021-10-05 14:35:33.681656"
$dataSheet.Range("F9").Value = "2021-10-05 14:35:33.681658"
$dataSheet.Range("F10").Value = "2021-10-05 14:35:33.681661"
$dataSheet.Range("F11").Value = "2021-10-05 14:35:33.681663"
$dataSheet.Range("F12").Value = "2021-10-05 14:35:33.681666"
$dataSheet.Range("F13").Value = "2021-10-05 14:35:33.681668"
$dataSheet.Range("F14").Value = "2021-10-05 14:35:33.681671"
$dataSheet.Range("F15").Value = "2021-10-05 14:35:33.681673"
$dataSheet.Range("F16").Value = "2021-10-05 14:35:33.681675"
$dataSheet.Range("F17").Value = "2021-10-05 14:35:33.681678"
$dataSheet.Range("F18").Value = "2021-10-05 14:35:33.681681"
$dataSheet.Range("F19").Value = "2021-10-05 14:35:33.681683"
$dataSheet.Range("F20").Value = "2021-10-05 14:35:33.681686"
$dataSheet.Range("F21").Value = "2021-10-05 14:35:33.681688"
$dataSheet.Range("F22").Value = "2021-10-05 14:35:33.681691"
$dataSheet.Range("F23").Value = "2021-10-05 14:35:33.681693"
$dataSheet.Range("F24").Value = "2021-10-05 14:35:33.681696"
$dataSheet.Range("F25").Value = "2021-10-05 14:35:33.681698"
$dataSheet.Range("F26").Value = "2021-10-05 14:35:33.681701"
$dataSheet.Range("F27").Value = "2021-10-05 14:35:33.681703"
$dataSheet.Range("F28").Value = "2021-10-05 14:35:33.681706"
$dataSheet.Range("F29").Value = "2021-10-05 14:35:33.681708"
$dataSheet.Range("F30").Value = "2021-10-05 14:35:33.681711"
$dataSheet.Range("F31").Value = "2021-10-05 14:35:33.681713"
$dataSheet.Range("F32").Value = "2021-10-05 14:35:33.681716"
$dataSheet.Range("F33").Value = "2021-10-05 14:35:33.681718"
$dataSheet.Range("F34").Value = "2021-10-05 14:35:33.681721"
$dataSheet.Range("F35").Value = "2021-10-05 14:35:33.681724"
$dataSheet.Range("F36").Value = "2021-10-05 14:35:33.681726"
$dataSheet.Range("F37").Value = "2021-10-05 14:35:33.681728"
$dataSheet.Range("F38").Value = "2021-10-05 14:35:33.681731"
$dataSheet.Range("F39").Value = "2021-10-05 14:35:33.681733"
$dataSheet.Range("F40").Value = "2021-10-05 14:35:33.681736"
$dataSheet.Range("F41").Value = "2021-10-05 14:35:33.681738"
$dataSheet.Range("F42").Value = "2021-10-05 14:35:33.681741"
$dataSheet.Range("F43").Value = "2021-10-05 14:35:33.681744"
$dataSheet.Range("F44").Value = "2021-10-05 14:35:33.681746"
$dataSheet.Range("F45").Value = "2021-10-05 14:35:33.681749"
$dataSheet.Range("F46").Value = "2021-10-05 14:35:33.681751"
$dataSheet.Range("F47").Value = "2021-10-05 14:35:33.681754"
$dataSheet.Range("F48").Value = "2021-10-05 14:35:33.681756"
$dataSheet.Range("F49").Value = "2021-10-05 14:35:33.681759"
$dataSheet.Range("F50").Value = "2021-10-05 14:35:33.681761"
$dataSheet.Range("F51").Value = "2021-10-05 14:35:33.681763"
$dataSheet.Range("F52").Value = "2021-10-05 14:35:33.681766"
$dataSheet.Range("F53").Value = "2021-10-05 14:35:33.681768"
$dataSheet.Range("F54").Value = "2021-10-05 14:35:33.681771"
$dataSheet.Range("F55").Value = "2021-10-05 14:35:33.681774"
$dataSheet.Range("F56").Value = "2021-10-05 14:35:33.681776"
$dataSheet.Range("F57").Value = "2021-10-05 14:35:33.681779"
$dataSheet.Range("F58").Value = "2021-10-05 14:35:33.681781"
$dataSheet.Range("F59").Value = "2021-10-05 14:35:33.681783"
$dataSheet.Range("F60").Value = "2021-10-05 14:35:33.681786"
$dataSheet.Range("F61").Value = "2021-10-05 14:35:33.681788"
$dataSheet.Range("F62").Value = "2021-10-05 14:35:33.681791"
$dataSheet.Range("F63").Value = "2021-10-05 14:35:33.681793"
$dataSheet.Range("F64").Value = "2021-10-05 14:35:33.681796"
$dataSheet.Range("F65").Value = "2021-10-05 14:35:33.681798"
$dataSheet.Range("F66").Value = "2021-10-05 14:35:33.681802"
$dataSheet.Range("F67").Value = "2021-10-05 14:35:33.681804"
$dataSheet.Range("F68").Value = "2021-10-05 14:35:33.681807"
$dataSheet.Range("F69").Value = "2021-10-05 14:35:33.681809"
$dataSheet.Range("F70").Value = "2021-10-05 14:35:33.681812"
$dataSheet.Range("F71").Value = "2021-10-05 14:35:33.681814"
$dataSheet.Range("F72").Value = "2021-10-05 14:35:33.681817"
$dataSheet.Range("F73").Value = "2021-10-05 14:35:33.681819"
$dataSheet.Range("F74").Value = "2021-10-05 14:35:33.681822"
$dataSheet.Range("F75").Value = "2021-10-05 14:35:33.681824"
$dataSheet.Range("F76").Value = "2021-10-05 14:35:33.681827"
$dataSheet.Range("F77").Value = "2021-10-05 14:35:33.681830"
$dataSheet.Range("F78").Value = "2021-10-05 14:35:33.681834"
$dataSheet.Range("F79").Value = "2021-10-05 14:35:33.681837"
$dataSheet.Range("F80").Value = "2021-10-05 14:35:33.681840"
$dataSheet.Range("F81").Value = "2021-10-05 14:35:33.681842"
$dataSheet.Range("F82").Value = "2021-10-05 14:35:33.681845"
$dataSheet.Range("F83").Value = "2021-10-05 14:35:33.681847"
$dataSheet.Range("F84").Value = "2021-10-05 14:35:33.681850"
$dataSheet.Range("F85").Value = "2021-10-05 14:35:33.681853"
$dataSheet.Range("F86").Value = "2021-10-05 14:35:33.681855"
$dataSheet.Range("F87").Value = "2021-10-05 14:35:33.681858"
$dataSheet.Range("F88").Value = "2021-10-05 14:35:33.681861"
$dataSheet.Range("F89").Value = "2021-10-05 14:35:33.681863"
$dataSheet.Range("F90").Value = "2021-10-05 14:35:33.681865"
$dataSheet.Range("F91").Value = "2021-10-05 14:35:33.681868"
$dataSheet.Range("F92").Value = "2021-10-05 14:35:33.681870"
$dataSheet.Range("F93").Value = "2021-10-05 14:35:33.681873"
$dataSheet.Range("F94").Value = "2021-10-05 14:35:33.681876"
$dataSheet.Range("F95").Value = "2021-10-05 14:35:33.681879"

# Add a new "metadata" worksheet positioned right after "data".
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1).
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row (A2:G2).
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Renal Ciliopathies and Nephronophthisis"
$metaSheet.Range("C2").Value = 193

# Force D2 ("1.2") to stay text instead of being parsed as a number, then
# strip the number-format style so it falls back to the default style (no
# explicit s= attribute), matching the source file.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.2"
$metaSheet.Range("D2").ClearFormats()

$metaSheet.Range("E2").Value = "2021-08-10T22:11:48.659778Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:33.678123"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/193/?format=json"

# Copy the header/index styling used on the "data" sheet (bold, bordered,
# centered) onto the new header row and the A2 index cell.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match sheetPr/outlinePr + page margins from the sibling "data" sheet.
$metaSheet.Outline.SummaryRow = 1
$metaSheet.Outline.SummaryColumn = 1

$metaSheet.PageSetup.LeftMargin = 54
$metaSheet.PageSetup.RightMargin = 54
$metaSheet.PageSetup.TopMargin = 72
$metaSheet.PageSetup.BottomMargin = 72
$metaSheet.PageSetup.HeaderMargin = 36
$metaSheet.PageSetup.FooterMargin = 36

$metaSheet.Range("A1").Select() | Out-Null

# Restore "data" as the active sheet/tab (unchanged by the source diff).
$dataSheet.Activate()
$dataSheet.Range("A1").Select() | Out-Null
